# edit.ps1 — apply the diff:
#  1) Re-style the three tables (slides 14, 15, 16) from the local custom
#     "Table_0" style {2CDD0541-BBCF-4698-B9B7-8AB1BAD5CEA5} to the built-in
#     table style {62DA24C5-A67C-469E-BE55-C173DB1A0E1B}.
#  2) Swap the presentation's active theme color scheme ("Integral" / Red
#     Violet) for the stock "Office" palette that previously only lived in
#     the notes-master theme part.

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------------
$newStyleId = "{62DA24C5-A67C-469E-BE55-C173DB1A0E1B}"
$tableSlides = @(14, 15, 16)
foreach ($idx in $tableSlides) {
    $slide = $p.Slides.Item($idx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2) Theme color scheme ---------------------------------------------
# RGB() builds the little-endian COM color value (0x00BBGGRR) PowerPoint
# expects for ColorFormat/ThemeColor RGB properties.
function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$newColors = @(
    (RGB 0x00 0x00 0x00),  # 1  dk1
    (RGB 0xFF 0xFF 0xFF),  # 2  lt1
    (RGB 0x44 0x54 0x6A),  # 3  dk2
    (RGB 0xE7 0xE6 0xE6),  # 4  lt2
    (RGB 0x5B 0x9B 0xD5),  # 5  accent1
    (RGB 0xED 0x7D 0x31),  # 6  accent2
    (RGB 0xA5 0xA5 0xA5),  # 7  accent3
    (RGB 0xFF 0xC0 0x00),  # 8  accent4
    (RGB 0x44 0x72 0xC4),  # 9  accent5
    (RGB 0x70 0xAD 0x47),  # 10 accent6
    (RGB 0x05 0x63 0xC1),  # 11 hlink
    (RGB 0x95 0x4F 0x72)   # 12 folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $newColors.Count; $i++) {
    $tcs.Colors($i).RGB = $newColors[$i - 1]
}
